$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "J"
$ws.Range("G8").Select()
